$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.136.25"
$ws.Range("E2").Value = "  +5.75%  "
$ws.Range("D3").Value = "2.278.06"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'232.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").Value = "'64.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.92%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.436"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.02%  "
$ws.Range("E10").Value = "  +17.35%  "
$ws.Range("D11").Value = "'56.67"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'25.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.40%  "
$ws.Range("D13").Value = "'0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "2.613.41"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "'15.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("D16").Value = "'5.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.90%  "
$ws.Range("D17").Value = "'0.831"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("D18").Value = "2.264.67"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "44.019.70"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.80%  "
$ws.Range("D21").Value = "'73.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "'258.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.41%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +6.09%  "
$ws.Range("D27").Value = "'10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("D28").Value = "'171.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("D29").Value = "'21.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.96%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").Value = "'2.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.06%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "'0.0685"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("D35").Value = "'4.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").Value = "'5.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'3.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.75%  "
$ws.Range("D38").Value = "'6.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.32%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").Value = "'0.0250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'8.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").Value = "'17.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.87%  "
$ws.Range("D44").Value = "'0.0972"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'4.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("D46").Value = "'98.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("B48").Value = "TerraClassic"
$ws.Range("C48").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D48").Value = "'0.000212"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.72%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.471.44"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.45%  "
$ws.Range("E51").Value = "  +1.94%  "
